# Trade #20 closed at 2026-02-17 04:08:13 - unknown UNKNOWN +0.000%
#
# - Summary sheet: Total Trades 19 -> 20, Win Rate % 31.58 -> 30
# - Strategy Status sheet: MarketMaking Trades 19 -> 20, Win Rate % 31.58 -> 30
# - All Trades / MarketMaking sheets: append trade #20 as new row 21

$wb = $excel.ActiveWorkbook

# --- Summary ---------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 20
$summary.Range("B9").Value = 30

# --- Strategy Status ---------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 20
$status.Range("G4").Value = 30

# --- Append the new closed trade row to both trade-log sheets ---------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 21

    $ws.Cells.Item($row, 1).Value = 20

    # Keep the date as literal text (matches the existing "Date" column
    # cells), not an auto-converted Excel date serial.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "04:08:07"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.5600000000000001
    $ws.Cells.Item($row, 7).Value = 0.5600000000000001
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.12
}
